$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 10.5
$ws.Range("C4").Value = 1.25
$ws.Range("C5").Value = 18

$ws.Range("C5").Select()
